$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Helper: locate a paragraph's 1-based index by its exact text (ignoring the
# trailing paragraph-mark character Word appends to Range.Text).
function Find-ParaIndex($doc, $text) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $idx
        }
    }
    return -1
}

# --- Part 1: paragraphs "    )", "  )" and "Fecha bien" -------------------
# Strip the en-US language formatting / proofErr markers and merge the
# "Fecha"/" bien" runs into a single plain run.

$i14 = Find-ParaIndex $d "    )"
$xml14 = "<w:p $wns><w:r><w:t xml:space='preserve'>    </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>"
$d.Paragraphs.Item($i14).Range.InsertXML($xml14)

$i15 = Find-ParaIndex $d "  )"
$xml15 = "<w:p $wns><w:r><w:t xml:space='preserve'>  )</w:t></w:r></w:p>"
$d.Paragraphs.Item($i15).Range.InsertXML($xml15)

$i16 = Find-ParaIndex $d "Fecha bien"
$xml16 = "<w:p $wns><w:r><w:t>Fecha bien</w:t></w:r></w:p>"
$d.Paragraphs.Item($i16).Range.InsertXML($xml16)

# --- Part 2: final paragraph "gravity " -----------------------------------
# Replace the truncated "gravity " paragraph with the full restored block of
# R code, keeping the trailing _GoBack bookmark on the very last paragraph.

$i65 = Find-ParaIndex $d "gravity "
$xml2 = @"
<w:p $wns>
  <w:pPr>
    <w:pBdr>
      <w:bottom w:val="thinThickThinMediumGap" w:sz="18" w:space="1" w:color="auto"/>
    </w:pBdr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">gravity </w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>source(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>here(</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">  "01_Scripts",</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">  "03_Resultados",</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">  </w:t>
  </w:r>
  <w:r>
    <w:t>"03_3_Analisis de errores",</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">  "013_Cálculo de </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>errores.R</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t>"</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:t>))</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@
$d.Paragraphs.Item($i65).Range.InsertXML($xml2)
